# Add files via upload
# - Insert three new columns (D:F) on the "axes" sheet, duplicating the
#   A/B/C (An/Ab/Or) header+value pair under new "*_arrow" headers, pushing
#   the old "Title" column (was D) out to G.
# - Make "axes" the active sheet/tab with selection on E5 (matches the
#   saved sheetView state in the target workbook).

$wb = $excel.ActiveWorkbook

$plys = $wb.Worksheets.Item("plys")
$axes = $wb.Worksheets.Item("axes")

# Insert 3 new columns before the existing "Title" column (D), shifting it to G.
[void]$axes.Columns("D:F").Insert()

# New header row (row 1)
$axes.Cells.Item(1, 4).Value = "A_arrow"
$axes.Cells.Item(1, 5).Value = "B_arrow"
$axes.Cells.Item(1, 6).Value = "C_arrow"

# New value row (row 2) - duplicate the An/Ab/Or values from columns A:C
$axes.Cells.Item(2, 4).Value = $axes.Cells.Item(2, 1).Value2
$axes.Cells.Item(2, 5).Value = $axes.Cells.Item(2, 2).Value2
$axes.Cells.Item(2, 6).Value = $axes.Cells.Item(2, 3).Value2

# Approximate the on-disk column width for the new columns (D:F).
$axes.Columns("D:F").ColumnWidth = 13.3

# Restore the original selection on "plys" (first sheet) so it is preserved
# once it stops being the active tab.
[void]$plys.Range("A9:XFD9").Select()

# Make "axes" the active sheet with the selection left on E5, matching the
# saved workbook state.
[void]$axes.Activate()
[void]$axes.Range("E5").Select()
